$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "F,S" cells to "F,S,M" (rows 4, 7, 10 - columns E, F, G)
$ws.Range("E4").Value = "F,S,M"
$ws.Range("F4").Value = "F,S,M"
$ws.Range("G4").Value = "F,S,M"

$ws.Range("E7").Value = "F,S,M"
$ws.Range("F7").Value = "F,S,M"
$ws.Range("G7").Value = "F,S,M"

$ws.Range("E10").Value = "F,S,M"
$ws.Range("F10").Value = "F,S,M"
$ws.Range("G10").Value = "F,S,M"

# Update the worksheet title cell (C2): "...elements and components" -> "...elements and element components"
$ws.Range("C2").Value = "Thinking about three gender-grouping Leslie matrix elements and element components"
